$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Safe (non-numeric-looking) text updates - direct assignment
$ws.Range('D2').Value = '62.674.33'
$ws.Range('E2').Value = '  +2.83%  '
$ws.Range('D3').Value = '2.943.94'
$ws.Range('E3').Value = '  +2.02%  '
$ws.Range('E4').Value = '  -0.03%  '
$ws.Range('E5').Value = '  +0.35%  '
$ws.Range('E6').Value = '  +7.13%  '
$ws.Range('D8').Value = '2.943.23'
$ws.Range('E8').Value = '  +2.01%  '
$ws.Range('E9').Value = '  +3.00%  '
$ws.Range('E10').Value = '  +1.77%  '
$ws.Range('E11').Value = '  +9.86%  '
$ws.Range('E12').Value = '  +2.39%  '
$ws.Range('E13').Value = '  +9.06%  '
$ws.Range('E14').Value = '  +0.76%  '
$ws.Range('E15').Value = '  -0.31%  '
$ws.Range('D16').Value = '3.433.11'
$ws.Range('E16').Value = '  +2.03%  '
$ws.Range('D17').Value = '62.636.39'
$ws.Range('E17').Value = '  +2.90%  '
$ws.Range('E18').Value = '  +2.56%  '
$ws.Range('D19').Value = '2.943.76'
$ws.Range('E19').Value = '  +1.94%  '
$ws.Range('E20').Value = '  +3.02%  '
$ws.Range('E21').Value = '  +2.29%  '
$ws.Range('E22').Value = '  +1.70%  '
$ws.Range('E23').Value = '  +1.01%  '
$ws.Range('E24').Value = '  +8.06%  '
$ws.Range('E25').Value = '  +0.75%  '
$ws.Range('E26').Value = '  +4.80%  '
$ws.Range('E27').Value = '  +2.79%  '
$ws.Range('E28').Value = '  -0.05%  '
$ws.Range('E29').Value = '  +10.94%  '
$ws.Range('E30').Value = '  +23.22%  '
$ws.Range('E31').Value = '  +2.46%  '
$ws.Range('E32').Value = '  +4.74%  '
$ws.Range('E33').Value = '  +5.87%  '
$ws.Range('E34').Value = '  +2.48%  '
$ws.Range('E35').Value = '  +0.08%  '
$ws.Range('E36').Value = '  +1.80%  '
$ws.Range('E37').Value = '  +11.62%  '
$ws.Range('E38').Value = '  +2.91%  '
$ws.Range('E39').Value = '  +1.52%  '
$ws.Range('E40').Value = '  +7.01%  '
$ws.Range('E41').Value = '  +1.13%  '
$ws.Range('E42').Value = '  +1.02%  '
$ws.Range('E43').Value = '  +5.00%  '
$ws.Range('E44').Value = '  +3.76%  '
$ws.Range('B45').Value = 'Monero'
$ws.Range('C45').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('E45').Value = '  +3.91%  '
$ws.Range('B46').Value = 'Maker'
$ws.Range('C46').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D46').Value = '2.705.53'
$ws.Range('E46').Value = '  +1.79%  '
$ws.Range('E47').Value = '  +3.93%  '
$ws.Range('E48').Value = '  +1.41%  '
$ws.Range('E50').Value = '  +2.51%  '
$ws.Range('E51').Value = '  +2.38%  '

# Numeric-looking values must be forced to remain text (matching original inline string formatting)
$numericTextValues = @{
    'D4' = '1.00'
    'D5' = '590.65'
    'D6' = '148.49'
    'D9' = '0.507'
    'D10' = '7.14'
    'D12' = '0.437'
    'D14' = '32.45'
    'D18' = '6.65'
    'D20' = '437.88'
    'D21' = '13.48'
    'D22' = '0.665'
    'D24' = '11.22'
    'D25' = '80.40'
    'D26' = '11.91'
    'D27' = '2.12'
    'D29' = '7.33'
    'D30' = '0.0000103'
    'D32' = '2.16'
    'D34' = '26.15'
    'D35' = '1.00'
    'D36' = '0.987'
    'D37' = '3.12'
    'D38' = '5.58'
    'D39' = '49.65'
    'D41' = '8.40'
    'D42' = '0.116'
    'D44' = '39.72'
    'D45' = '135.95'
    'D48' = '357.32'
}
foreach ($cell in $numericTextValues.Keys) {
    $ws.Range($cell).NumberFormat = "@"
    $ws.Range($cell).Value = $numericTextValues[$cell]
    $ws.Range($cell).Style = "Normal"
}
